$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J2").Value = "Test119/1"
$ws.Range("M2").Value = "'3000"
$ws.Range("N2").Value = "D200E"
$ws.Range("O2").Value = "'60"
Write-Host "done"
